$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8019073605537415
$ws.Range("B1").Value = 0.8509652614593506
$ws.Range("C1").Value = 5.025641441345215
$ws.Range("D1").Value = 1.793330311775208
$ws.Range("E1").Value = 0.7575253248214722
